$wb = $excel.ActiveWorkbook

# --- Rename sheets ---------------------------------------------------------
$wsMain = $wb.Worksheets.Item(1)
$wsMain.Name = "Export as TSV"

$wsUnit = $wb.Worksheets.Item(5)
$wsUnit.Name = "bulk_rna_yield...ssue_unit list"

# --- Freeze header row on the main sheet -----------------------------------
$wsMain.Activate()
$wsMain.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- Data validation error messages / titles --------------------------------
$dv = $wsMain.Range("I2:I1048576").Validation
$dv.ErrorTitle = "Value must come from list"
$dv.ErrorMessage = "Value must be one of: sequence."

$dv = $wsMain.Range("J2:J1048576").Validation
$dv.ErrorTitle = "Value must come from list"
$dv.ErrorMessage = "Value must be one of: bulk RNA."

$dv = $wsMain.Range("K2:K1048576").Validation
$dv.ErrorTitle = "Value must come from list"
$dv.ErrorMessage = "Value must be one of: RNA."

$dv = $wsMain.Range("L2:L1048576").Validation
$dv.ErrorTitle = "Not a boolean"
$dv.ErrorMessage = 'The values in this column must be "TRUE" or "FALSE".'

$dv = $wsMain.Range("P2:P1048576").Validation
$dv.ErrorTitle = "Not a number"
$dv.ErrorMessage = "The values in this column must be numbers."

# Column Q also needs its source list formula repointed at the renamed sheet
$dv = $wsMain.Range("Q2:Q1048576").Validation
$dv.ErrorTitle = "Value must come from list"
$dv.ErrorMessage = "Value must be one of: ng/mg."
$dv.Formula1 = "='bulk_rna_yield...ssue_unit list'!`$A`$1:`$A`$1"

$dv = $wsMain.Range("R2:R1048576").Validation
$dv.ErrorTitle = "Not a number"
$dv.ErrorMessage = "The values in this column must be numbers."

$dv = $wsMain.Range("S2:S1048576").Validation
$dv.ErrorTitle = "Not a number"
$dv.ErrorMessage = "The values in this column must be numbers."

$dv = $wsMain.Range("T2:T1048576").Validation
$dv.ErrorTitle = "Value must come from list"
$dv.ErrorMessage = "Value must be one of: ug."

$dv = $wsMain.Range("W2:W1048576").Validation
$dv.ErrorTitle = "Value must come from list"
$dv.ErrorMessage = "Value must be one of: single-end / paired-end."

$dv = $wsMain.Range("Y2:Y1048576").Validation
$dv.ErrorTitle = "Not an integer"
$dv.ErrorMessage = "The values in this column must be integers."

$dv = $wsMain.Range("Z2:Z1048576").Validation
$dv.ErrorTitle = "Not a number"
$dv.ErrorMessage = "The values in this column must be numbers."

$dv = $wsMain.Range("AA2:AA1048576").Validation
$dv.ErrorTitle = "Value must come from list"
$dv.ErrorMessage = "Value must be one of: ng."

$dv = $wsMain.Range("AB2:AB1048576").Validation
$dv.ErrorTitle = "Not an integer"
$dv.ErrorMessage = "The values in this column must be integers."

$dv = $wsMain.Range("AE2:AE1048576").Validation
$dv.ErrorTitle = "Not a number"
$dv.ErrorMessage = "The values in this column must be numbers."

$dv = $wsMain.Range("AF2:AF1048576").Validation
$dv.ErrorTitle = "Not a number"
$dv.ErrorMessage = "The values in this column must be numbers."

Write-Output "done"
